# Weekly update: insert a new price record (week of 2021-10-05) as row 312
# in the "Ajo" (Chino / Primera) series, pushing the existing rows 312-351
# down to 313-352.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 312, shifting rows 312:351 -> 313:352.
$ws.Rows.Item(312).Insert()

# Populate the newly inserted row 312 with the new weekly record.
$ws.Cells.Item(312, 1).Value = 6
$ws.Cells.Item(312, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(312, 3).Value = "Metropolitana"
$ws.Cells.Item(312, 4).Value = 44474
$ws.Cells.Item(312, 5).Value = 13
$ws.Cells.Item(312, 6).Value = 100112003
$ws.Cells.Item(312, 7).Value = "Ajo"
$ws.Cells.Item(312, 8).Value = "Chino"
$ws.Cells.Item(312, 9).Value = "Primera"
$ws.Cells.Item(312, 10).Value = 2800
$ws.Cells.Item(312, 11).Value = 14500
$ws.Cells.Item(312, 12).Value = 15000
$ws.Cells.Item(312, 13).Value = 14732
$ws.Cells.Item(312, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(312, 15).Value = "China"
$ws.Cells.Item(312, 16).Value = 1473
$ws.Cells.Item(312, 17).Value = 10
$ws.Cells.Item(312, 18).Value = "Hortaliza"
